# Insert a new weekly record for "Femacal de La Calera - Apio" (Hortaliza)
# just above the current row 208, pushing all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 208 (and everything below it) down by one row.
$ws.Rows.Item(208).Insert()

# Populate the newly inserted row 208 with the new weekly data point.
$ws.Range("A208").Value = 3
$ws.Range("B208").Value = "Femacal de La Calera"
$ws.Range("C208").Value = "Coquimbo"
$ws.Range("D208").Value = 44474
$ws.Range("E208").Value = 5
$ws.Range("F208").Value = 100112017
$ws.Range("G208").Value = "Apio"
$ws.Range("H208").Value = "Americana (o)"
$ws.Range("I208").Value = "Primera"
$ws.Range("J208").Value = 160
$ws.Range("K208").Value = 9000
$ws.Range("L208").Value = 9000
$ws.Range("M208").Value = 9000
$ws.Range("N208").Value = "$/docena de matas"
$ws.Range("O208").Value = "Pan de Azúcar"
$ws.Range("P208").Value = 1500
$ws.Range("Q208").Value = 6
$ws.Range("R208").Value = "Hortaliza"
